$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add the new "Drums" row of data in row 21 (previously blank except A21 styling)
$ws.Range("A21").Value = "drums"
$ws.Range("B21").Value = "Drums"
$ws.Range("C21").Value = "NA"
$ws.Range("D21").Value = "Schlagzeug"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = "auto"
$ws.Range("J21").Value = "NA"

# Update the active selection to reflect where the user ended up after entering data
$ws.Range("D24").Select()
